# This script updates the LR-pair TPM statistics sheet to match the newly
# computed values ("update scripts wuth new tpm"). Column D (Target cluster)
# keeps the same cluster-name text per row; only the underlying shared-string
# bookkeeping differs in the source diff, so we (re)assert the same label text
# to make that explicit, and then update every changed numeric metric E:T.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4913865
$ws.Range("H2").Value = 0.982773
$ws.Range("I2").Value = 0.7978605686549953
$ws.Range("J2").Value = 0.789798398179267
$ws.Range("M2").Value = 10.306905
$ws.Range("N2").Value = 20.61381
$ws.Range("O2").Value = 0.04169074224953703
$ws.Range("P2").Value = 0.0286470056427464
$ws.Range("Q2").Value = 5.0646739737825
$ws.Range("R2").Value = 20.25869589513
$ws.Range("S2").Value = 0.03326339931886445
$ws.Range("T2").Value = 0.02262535916927353

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4913865
$ws.Range("H3").Value = 0.982773
$ws.Range("I3").Value = 0.7978605686549953
$ws.Range("J3").Value = 0.789798398179267
$ws.Range("O3").Value = 0.04663503533846117
$ws.Range("P3").Value = 0.04806657479834425
$ws.Range("Q3").Value = 5.665316494761
$ws.Range("R3").Value = 33.991898968566
$ws.Range("S3").Value = 0.03720825581439043
$ws.Range("T3").Value = 0.03796290378169621

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4913865
$ws.Range("H4").Value = 0.982773
$ws.Range("I4").Value = 0.7978605686549953
$ws.Range("J4").Value = 0.789798398179267
$ws.Range("M4").Value = 65.286547
$ws.Range("N4").Value = 195.859641
$ws.Range("O4").Value = 0.2640797216370273
$ws.Range("P4").Value = 0.2721860850038534
$ws.Range("Q4").Value = 32.0809278274155
$ws.Range("R4").Value = 192.485566964493
$ws.Range("S4").Value = 0.2106987968755715
$ws.Range("T4").Value = 0.2149721339427293

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4913865
$ws.Range("H5").Value = 0.982773
$ws.Range("I5").Value = 0.7978605686549953
$ws.Range("J5").Value = 0.789798398179267
$ws.Range("M5").Value = 11.7817895
$ws.Range("N5").Value = 23.563579
$ws.Range("O5").Value = 0.04765655153344304
$ws.Range("P5").Value = 0.03274629874711665
$ws.Range("Q5").Value = 5.78941230614175
$ws.Range("R5").Value = 23.157649224567
$ws.Range("S5").Value = 0.03802328330660895
$ws.Range("T5").Value = 0.02586297429677246

# Row 6
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4913865
$ws.Range("H6").Value = 0.982773
$ws.Range("I6").Value = 0.7978605686549953
$ws.Range("J6").Value = 0.789798398179267
$ws.Range("M6").Value = 87.427925
$ws.Range("N6").Value = 262.283775
$ws.Range("O6").Value = 0.3536401166583814
$ws.Range("P6").Value = 0.3644956843216187
$ws.Range("Q6").Value = 42.9609020680125
$ws.Range("R6").Value = 257.765412408075
$ws.Range("S6").Value = 0.2821555045762751
$ws.Range("T6").Value = 0.2878781076204702

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.4913865
$ws.Range("H7").Value = 0.982773
$ws.Range("I7").Value = 0.7978605686549953
$ws.Range("J7").Value = 0.789798398179267
$ws.Range("M7").Value = 60.89045733333333
$ws.Range("N7").Value = 182.671372
$ws.Range("O7").Value = 0.24629783258315
$ws.Range("P7").Value = 0.2538583514863204
$ws.Range("Q7").Value = 29.920748712426
$ws.Range("R7").Value = 179.524492274556
$ws.Range("S7").Value = 0.1965113287632849
$ws.Range("T7").Value = 0.2004969193683252

# Row 8
$ws.Range("I8").Value = 0.02041576811073366
$ws.Range("J8").Value = 0.03031420824312443
$ws.Range("M8").Value = 10.306905
$ws.Range("N8").Value = 20.61381
$ws.Range("O8").Value = 0.04169074224953703
$ws.Range("P8").Value = 0.0286470056427464
$ws.Range("Q8").Value = 0.129595587835
$ws.Range("R8").Value = 0.7775735270099999
$ws.Range("S8").Value = 0.0008511485261309146
$ws.Range("T8").Value = 0.0008684112945961752

# Row 9
$ws.Range("I9").Value = 0.02041576811073366
$ws.Range("J9").Value = 0.03031420824312443
$ws.Range("O9").Value = 0.04663503533846117
$ws.Range("P9").Value = 0.04806657479834425
$ws.Range("S9").Value = 0.0009520900673058928
$ws.Range("T9").Value = 0.001457100157970724

# Row 10
$ws.Range("I10").Value = 0.02041576811073366
$ws.Range("J10").Value = 0.03031420824312443
$ws.Range("M10").Value = 65.286547
$ws.Range("N10").Value = 195.859641
$ws.Range("O10").Value = 0.2640797216370273
$ws.Range("P10").Value = 0.2721860850038534
$ws.Range("Q10").Value = 0.8208912797956666
$ws.Range("R10").Value = 7.388021518160999
$ws.Range("S10").Value = 0.005391390359688644
$ws.Range("T10").Value = 0.008251105661687581

# Row 11
$ws.Range("D11").Value = "MuSCs"
$ws.Range("I11").Value = 0.02041576811073366
$ws.Range("J11").Value = 0.03031420824312443
$ws.Range("M11").Value = 11.7817895
$ws.Range("N11").Value = 23.563579
$ws.Range("O11").Value = 0.04765655153344304
$ws.Range("P11").Value = 0.03274629874711665
$ws.Range("Q11").Value = 0.1481402939098333
$ws.Range("R11").Value = 0.8888417634589998
$ws.Range("S11").Value = 0.0009729451050640017
$ws.Range("T11").Value = 0.0009926781194116588

# Row 12
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("I12").Value = 0.02041576811073366
$ws.Range("J12").Value = 0.03031420824312443
$ws.Range("M12").Value = 87.427925
$ws.Range("N12").Value = 262.283775
$ws.Range("O12").Value = 0.3536401166583814
$ws.Range("P12").Value = 0.3644956843216187
$ws.Range("Q12").Value = 1.099289586308333
$ws.Range("R12").Value = 9.893606276774999
$ws.Range("S12").Value = 0.007219834616350314
$ws.Range("T12").Value = 0.01104939807824569

# Row 13
$ws.Range("I13").Value = 0.02041576811073366
$ws.Range("J13").Value = 0.03031420824312443
$ws.Range("M13").Value = 60.89045733333333
$ws.Range("N13").Value = 182.671372
$ws.Range("O13").Value = 0.24629783258315
$ws.Range("P13").Value = 0.2538583514863204
$ws.Range("Q13").Value = 0.7656163136902221
$ws.Range("R13").Value = 6.890546823212
$ws.Range("S13").Value = 0.005028359436193892
$ws.Range("T13").Value = 0.007695514931212594

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.5
$ws.Range("G14").Value = 0.11192
$ws.Range("H14").Value = 0.22384
$ws.Range("I14").Value = 0.181723663234271
$ws.Range("J14").Value = 0.1798873935776086
$ws.Range("M14").Value = 10.306905
$ws.Range("N14").Value = 20.61381
$ws.Range("O14").Value = 0.04169074224953703
$ws.Range("P14").Value = 0.0286470056427464
$ws.Range("Q14").Value = 1.1535488076
$ws.Range("R14").Value = 4.6141952304
$ws.Range("S14").Value = 0.00757619440454166
$ws.Range("T14").Value = 0.005153235178876696

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.5
$ws.Range("G15").Value = 0.11192
$ws.Range("H15").Value = 0.22384
$ws.Range("I15").Value = 0.181723663234271
$ws.Range("J15").Value = 0.1798873935776086
$ws.Range("O15").Value = 0.04663503533846117
$ws.Range("P15").Value = 0.04806657479834425
$ws.Range("Q15").Value = 1.290353361546667
$ws.Range("R15").Value = 7.74212016928
$ws.Range("S15").Value = 0.008474689456764843
$ws.Range("T15").Value = 0.008646570858677314

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.5
$ws.Range("G16").Value = 0.11192
$ws.Range("H16").Value = 0.22384
$ws.Range("I16").Value = 0.181723663234271
$ws.Range("J16").Value = 0.1798873935776086
$ws.Range("M16").Value = 65.286547
$ws.Range("N16").Value = 195.859641
$ws.Range("O16").Value = 0.2640797216370273
$ws.Range("P16").Value = 0.2721860850038534
$ws.Range("Q16").Value = 7.306870340240001
$ws.Range("R16").Value = 43.84122204144
$ws.Range("S16").Value = 0.04798953440176717
$ws.Range("T16").Value = 0.04896284539943661

# Row 17
$ws.Range("D17").Value = "MuSCs"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.5
$ws.Range("G17").Value = 0.11192
$ws.Range("H17").Value = 0.22384
$ws.Range("I17").Value = 0.181723663234271
$ws.Range("J17").Value = 0.1798873935776086
$ws.Range("M17").Value = 11.7817895
$ws.Range("N17").Value = 23.563579
$ws.Range("O17").Value = 0.04765655153344304
$ws.Range("P17").Value = 0.03274629874711665
$ws.Range("Q17").Value = 1.31861788084
$ws.Range("R17").Value = 5.27447152336
$ws.Range("S17").Value = 0.008660323121770083
$ws.Range("T17").Value = 0.005890646330932523

# Row 18
$ws.Range("D18").Value = "Neutrophils"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.5
$ws.Range("G18").Value = 0.11192
$ws.Range("H18").Value = 0.22384
$ws.Range("I18").Value = 0.181723663234271
$ws.Range("J18").Value = 0.1798873935776086
$ws.Range("M18").Value = 87.427925
$ws.Range("N18").Value = 262.283775
$ws.Range("O18").Value = 0.3536401166583814
$ws.Range("P18").Value = 0.3644956843216187
$ws.Range("Q18").Value = 9.784933366000001
$ws.Range("R18").Value = 58.709600196
$ws.Range("S18").Value = 0.064264777465756
$ws.Range("T18").Value = 0.0655681786229028

# Row 19
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.5
$ws.Range("G19").Value = 0.11192
$ws.Range("H19").Value = 0.22384
$ws.Range("I19").Value = 0.181723663234271
$ws.Range("J19").Value = 0.1798873935776086
$ws.Range("M19").Value = 60.89045733333333
$ws.Range("N19").Value = 182.671372
$ws.Range("O19").Value = 0.24629783258315
$ws.Range("P19").Value = 0.2538583514863204
$ws.Range("Q19").Value = 6.814859984746667
$ws.Range("R19").Value = 40.88915990848
$ws.Range("S19").Value = 0.04475814438367121
$ws.Range("T19").Value = 0.04566591718678262
